$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "ff"
$ws.Range("B7").Value = 12.6

$ws.Range("A8").Value = "gg"
$ws.Range("B8").Value = 14.3

$ws.Range("C8").Select()
